# Restore deployed sources and recovered app.html; prepare local rebuild
#
# 1) Update the "Weekly Timesheet" sheet: new hours/client/rate/total figures
#    for the week, plus a corrected subtotal row.
# 2) Add a new "Jason Schema" sheet: a flat, one-row-per-shift export keyed
#    by employee, mirroring the timesheet data in a wider, denormalised
#    layout used for downstream ingestion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# ---------------------------------------------------------------------------
# 1) Weekly Timesheet — refresh the daily rows and the subtotal
# ---------------------------------------------------------------------------

# Row 2: 2026-01-19 / Hall
$ws.Cells.Item(2, 2).Value = "Hall"
$ws.Cells.Item(2, 3).Value = 8
$ws.Cells.Item(2, 5).Value = 100
$ws.Cells.Item(2, 6).Value = 800

# Row 3: 2026-01-20 / McGill
$ws.Cells.Item(3, 2).Value = "McGill"
$ws.Cells.Item(3, 3).Value = 8
$ws.Cells.Item(3, 5).Value = 90
$ws.Cells.Item(3, 6).Value = 720

# Row 4: 2026-01-21 / Bryan
$ws.Cells.Item(4, 2).Value = "Bryan"
$ws.Cells.Item(4, 3).Value = 10
$ws.Cells.Item(4, 5).Value = 90
$ws.Cells.Item(4, 6).Value = 900

# Row 5: 2026-01-22 / McGill
$ws.Cells.Item(5, 2).Value = "McGill"
$ws.Cells.Item(5, 3).Value = 8
$ws.Cells.Item(5, 5).Value = 90
$ws.Cells.Item(5, 6).Value = 720

# Row 6: date moves from 2026-01-23 to 2026-01-25 / McGill
# (force text so it isn't auto-coerced into a date serial, matching the
# plain-text dates already used in column A)
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "2026-01-25"
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(6, 2).Value = "McGill"
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 5).Value = 90
$ws.Cells.Item(6, 6).Value = 720

# Row 8: SUBTOTAL — hours 38 -> 42, total 3730 -> 3860
$ws.Cells.Item(8, 3).Value = 42
$ws.Cells.Item(8, 4).Value = "Reg: 42 / OT: 0"
$ws.Cells.Item(8, 6).Value = 3860

# ---------------------------------------------------------------------------
# 2) Add the "Jason Schema" sheet right after "Weekly Timesheet"
# ---------------------------------------------------------------------------

$schemaWs = $wb.Worksheets.Add([Type]::Missing, $ws)
$schemaWs.Name = "Jason Schema"

# Column widths (character units) matching the target layout
$schemaWs.Columns.Item(1).ColumnWidth = 19.140625   # A Employee       -> 20
$schemaWs.Columns.Item(2).ColumnWidth = 17.1875     # B Employee ID    -> 18
$schemaWs.Columns.Item(3).ColumnWidth = 11.1328125  # C Date           -> 12
$schemaWs.Columns.Item(4).ColumnWidth = 24.21875    # D Client         -> 25
$schemaWs.Columns.Item(5).ColumnWidth = 7.2265625   # E Hours          -> 8
$schemaWs.Columns.Item(6).ColumnWidth = 9.1796875   # F Rate           -> 10
$schemaWs.Columns.Item(7).ColumnWidth = 11.1328125  # G Total          -> 12
$schemaWs.Columns.Item(8).ColumnWidth = 9.1796875   # H Type           -> 10
$schemaWs.Columns.Item(9).ColumnWidth = 29.1015625  # I Notes          -> 30

# Header row
$headers = @("Employee", "Employee ID", "Date", "Client", "Hours", "Rate", "Total", "Type", "Notes")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $schemaWs.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$schemaWs.Range("A1:I1").Font.Bold = $true
$schemaWs.Range("F1:G1").NumberFormat = '"$"#,##0.00'

# Data rows — one row per shift, employee columns repeated for each
$employee = "Chris Jacobi"
$employeeId = "emp_q3WGXgczT8gssfCO"

$rows = @(
    @("2026-01-19", "Hall",   8,  100, 800, "Regular"),
    @("2026-01-20", "McGill", 8,  90,  720, "Regular"),
    @("2026-01-21", "Bryan",  10, 90,  900, "Regular"),
    @("2026-01-22", "McGill", 8,  90,  720, "Regular"),
    @("2026-01-25", "McGill", 8,  90,  720, "Regular")
)

# Pre-format the Date column as text so the "2026-01-19"-style strings
# aren't auto-coerced into date serials when assigned below.
$schemaWs.Range("C2:C6").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $schemaWs.Cells.Item($excelRow, 1).Value = $employee
    $schemaWs.Cells.Item($excelRow, 2).Value = $employeeId
    $schemaWs.Cells.Item($excelRow, 3).Value = $row[0]
    $schemaWs.Cells.Item($excelRow, 4).Value = $row[1]
    $schemaWs.Cells.Item($excelRow, 5).Value = $row[2]
    $schemaWs.Cells.Item($excelRow, 6).Value = $row[3]
    $schemaWs.Cells.Item($excelRow, 7).Value = $row[4]
    $schemaWs.Cells.Item($excelRow, 8).Value = $row[5]
    $schemaWs.Cells.Item($excelRow, 9).Value = ""
}

$schemaWs.Range("C2:C6").Style = "Normal"
$schemaWs.Range("F2:G6").NumberFormat = '"$"#,##0.00'
